# Begonnen met Matlab files aan te passen voor scenarios
#
# Duplicate the "flights" sheet into two new scenario sheets
# ("Scenario1" / "Scenario2"), change their Type column (B) to a single
# runway-configuration value each, and append a couple of blank rows at
# the bottom of each, mirroring the manual "save as new scenario" steps
# the author performed in Excel.

$wb = $excel.ActiveWorkbook

$flights = $wb.Worksheets.Item("flights")

# --- Scenario1: copy of "flights", placed right after it ---------------
$flights.Copy($null, $flights)
$scenario1 = $wb.Worksheets.Item($wb.Worksheets.Count)
$scenario1.Name = "Scenario1"

# --- Scenario2: copy of "flights", placed right after Scenario1 --------
$flights.Copy($null, $scenario1)
$scenario2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$scenario2.Name = "Scenario2"

# Scenario1 uses runway configuration 3 for every flight
for ($r = 2; $r -le 11; $r++) {
    $scenario1.Cells.Item($r, 2).Value = 3
}

# Scenario2 uses runway configuration 2 for every flight
for ($r = 2; $r -le 11; $r++) {
    $scenario2.Cells.Item($r, 2).Value = 2
}

# Append two blank rows (12-13) below the table on Scenario1, and three
# (12-14) on Scenario2, matching the extra room left for future entries.
$scenario1.Range("A12:D13").Font.Bold = $false

$scenario2.Range("A12:D14").Font.Bold = $false

# --- Selections / active sheet, matching the final manual edit state ---
$flights.Activate()
$flights.Range("D3").Select() | Out-Null

$scenario1.Activate()
$scenario1.Range("C12").Select() | Out-Null

$scenario2.Activate()
$scenario2.Range("C12").Select() | Out-Null
